# Split the surname in each team member's line into its own run
# (mirrors Word auto-splitting a run when a spell-check-style
# property toggle is applied only to the surname substring) and
# split "drawio" out of the UML-diagrams sentence the same way.
$d = $word.ActiveDocument

$surnames = @("Serpatowska", "Szypulski", "J" + [char]0x0105 + "der", "" + [char]0x017B + "ywko", "drawio")

foreach ($name in $surnames) {
    $rng = $d.Content
    $found = $rng.Find.Execute($name, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Bold = 1
        $rng.Bold = 0
    }
}

# Remove the "Tests" section entirely: everything from the blank
# paragraph right after the UML-diagrams paragraph through the very
# last paragraph of the document ("Result:" of Test 4).
$umlRng = $d.Content
$null = $umlRng.Find.Execute("Due to the complexity")
$umlPara = $umlRng.Paragraphs(1)

$startPara = $umlPara.Next()
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$delRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
$delRange.Delete()
